# bot5 finalizado3 ruta logPath
#
# - Update the SAP GUI logon path stored on the "Rutas" sheet: it moves from
#   the old D:\...\ERPSAP\SAPgui\saplogon.exe install location to the new
#   C:\...\SAP\FrontEnd\SAPgui\saplogon.exe location.
# - "Rutas" becomes the active/selected sheet tab (previously
#   "parametrosInicio" was the active tab).

$wb = $excel.ActiveWorkbook

$wsRutas = $wb.Worksheets.Item("Rutas")

# New SAP GUI logon executable path
$wsRutas.Range("B2").Value = "C:\Program Files (x86)\SAP\FrontEnd\SAPgui\saplogon.exe"

# "Rutas" is now the selected/active sheet tab
$wsRutas.Activate()
